$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete row 7 first (LARISSA / 004363260 / 1045.41) so row indices for the
# earlier row stay valid, then delete row 5 (MARINA / 004556150 / 21998.48).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()
